# Bayes net figure: add a "Component Event" / composite-event slide.
#
# 1. Slide 1 & Slide 2: rename the "Attack Event" node label to "Event".
# 2. Duplicate Slide 2 to create a new Slide 3, then edit its copy of the
#    node label to "Component Event" and add a third, higher-level
#    "Composite Event" / "Hypothesis" node (ellipse + 2 connectors + label)
#    above it, wired into the existing diagram.

$p = $ppt.ActivePresentation

# --- Slide 1: "Attack Event" -> "Event" -----------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(5).TextFrame.TextRange.Paragraphs(1).Text = "Event"

# --- Duplicate Slide 2 (before editing its text) to build Slide 3 ---------
$s2 = $p.Slides.Item(2)
$s3 = $s2.Duplicate()

# --- Slide 2: "Attack Event" -> "Event" ------------------------------------
$s2.Shapes.Item(5).TextFrame.TextRange.Paragraphs(1).Text = "Event"

# --- Slide 3: relabel the copied node and add the composite-event node ----
$s3.Shapes.Item(5).TextFrame.TextRange.Paragraphs(1).Text = "Component Event"

# New "Composite Event" / "Hypothesis" label box (clone of the existing
# "Component Event" / "Hypotheses" label box, repositioned above it).
$lbl = $s3.Shapes.Item(5).Duplicate()
$lbl.Left = 235.72913385826772
$lbl.Top = 82.04866141732283
$lbl.Width = 128.94803149606298
$lbl.Height = 50.286614173228344
$lbl.TextFrame.TextRange.Paragraphs(1).Text = "Composite Event"
$lbl.TextFrame.TextRange.Paragraphs(2).Text = "Hypothesis"

# New ellipse node for the composite event (clone of an existing node).
$node = $s3.Shapes.Item(1).Duplicate()
$node.Left = 201.6
$node.Top = 101.52448818897638
$node.Width = 22.0251968503937
$node.Height = 24.094488188976378

# Connector from the composite-event node down into the diagram (flipped).
$conn1 = $s3.Shapes.Item(2).Duplicate()
$conn1.Left = 180.0
$conn1.Top = 122.22937007874016
$conn1.Width = 29.315905511811025
$conn1.Height = 28.970629921259842

# Second connector from the composite-event node down into the diagram.
$conn2 = $s3.Shapes.Item(4).Duplicate()
$conn2.Left = 209.31590551181102
$conn2.Top = 111.85566929133859
$conn2.Width = 35.484015748031496
$conn2.Height = 40.10976377952756
